$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "BetaFiberA"

# Add the new row of averaged-intensity data (row 16)
# Copy the style from the row above (A15) so the same formatting (bold,
# centered, bordered) is reused instead of creating a brand-new style.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.010501375714796
$ws.Range("D16").Value = 0.9862014279361111
$ws.Range("E16").Value = 0.9964998363637791
$ws.Range("F16").Value = 0.9879085223676277
$ws.Range("G16").Value = 1.010501375714796
$ws.Range("H16").Value = 0.9862014279361111
$ws.Range("I16").Value = 0.995867728876934
$ws.Range("J16").Value = 0.9923383171122281
$ws.Range("K16").Value = 0.9982352941176471
$ws.Range("L16").Value = 0.9805882352941176
$ws.Range("M16").Value = 1.010501375714796
$ws.Range("N16").Value = 0.9913506321499451
$ws.Range("O16").Value = 0.9952777905955785
$ws.Range("P16").Value = 0.9935175922229051
